$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.003.20"
$ws.Range("E2").Value = "  -3.63%  "

$ws.Range("D3").Value = "1.867.87"
$ws.Range("E3").Value = "  -2.88%  "

$ws.Range("D4").Value = "1.006"
$ws.Range("E4").Value = "  +0.13%  "

$ws.Range("D5").Value = "318.28"
$ws.Range("E5").Value = "  -2.26%  "

$ws.Range("E6").Value = "  +0.18%  "

$ws.Range("D7").Value = "0.4327"
$ws.Range("E7").Value = "  -5.88%  "

$ws.Range("D8").Value = "0.3707"
$ws.Range("E8").Value = "  -2.92%  "

$ws.Range("D9").Value = "0.07399"
$ws.Range("E9").Value = "  -4.57%  "

$ws.Range("D10").Value = "0.9308"
$ws.Range("E10").Value = "  -5.03%  "

$ws.Range("D11").Value = "21.18"
$ws.Range("E11").Value = "  -6.63%  "

$ws.Range("D12").Value = "1.900.29"
$ws.Range("E12").Value = "  -0.57%  "

$ws.Range("D13").Value = "6.713"
$ws.Range("E13").Value = "  -3.65%  "

$ws.Range("D14").Value = "5.418"
$ws.Range("E14").Value = "  -4.79%  "

$ws.Range("D15").Value = "0.06868"
$ws.Range("E15").Value = "  -2.33%  "

$ws.Range("E16").Value = "  +0.05%  "

$ws.Range("D17").Value = "80.02"
$ws.Range("E17").Value = "  -5.05%  "

$ws.Range("D18").Value = "0.000008971"
$ws.Range("E18").Value = "  -5.89%  "

$ws.Range("D19").Value = "1.005"
$ws.Range("E19").Value = "  +0.28%  "

$ws.Range("D20").Value = "15.72"
$ws.Range("E20").Value = "  -6.03%  "

$ws.Range("D21").Value = "28.008.81"
$ws.Range("E21").Value = "  -3.59%  "

$ws.Range("D22").Value = "5.098"
$ws.Range("E22").Value = "  -4.52%  "

$ws.Range("D23").Value = "10.97"
$ws.Range("E23").Value = "  -0.05%  "

$ws.Range("D24").Value = "2.202.66"
$ws.Range("E24").Value = "  +2.69%  "

$ws.Range("D25").Value = "2.049"
$ws.Range("E25").Value = "  -1.29%  "

$ws.Range("D26").Value = "154.11"
$ws.Range("E26").Value = "  -2.10%  "

$ws.Range("D27").Value = "18.47"
$ws.Range("E27").Value = "  -3.19%  "

$ws.Range("D28").Value = "5.454"
$ws.Range("E28").Value = "  -3.55%  "

$ws.Range("D29").Value = "112.86"
$ws.Range("E29").Value = "  -4.35%  "

$ws.Range("D30").Value = "1.680"
$ws.Range("E30").Value = "  -8.53%  "

$ws.Range("D31").Value = "0.08971"
$ws.Range("E31").Value = "  -4.01%  "

$ws.Range("D32").Value = "0.8037"
$ws.Range("E32").Value = "  -6.35%  "

$ws.Range("D33").Value = "4.746"
$ws.Range("E33").Value = "  -7.03%  "

$ws.Range("D34").Value = "1.169"
$ws.Range("E34").Value = "  -5.97%  "

$ws.Range("D35").Value = "2.956"
$ws.Range("E35").Value = "  -2.03%  "

$ws.Range("D36").Value = "1.005"
$ws.Range("E36").Value = "  +0.22%  "

$ws.Range("D37").Value = "0.05493"
$ws.Range("E37").Value = "  -3.34%  "

$ws.Range("E38").Value = "  -3.49%  "

$ws.Range("D39").Value = "0.01968"
$ws.Range("E39").Value = "  -3.84%  "

$ws.Range("D40").Value = "3.016"
$ws.Range("E40").Value = "  -2.85%  "

$ws.Range("D41").Value = "0.5224"
$ws.Range("E41").Value = "  -5.12%  "

$ws.Range("D42").Value = "6.986"
$ws.Range("E42").Value = "  -6.26%  "

$ws.Range("D43").Value = "0.1685"
$ws.Range("E43").Value = "  -3.99%  "

$ws.Range("D44").Value = "8.708"
$ws.Range("E44").Value = "  -7.26%  "

$ws.Range("D45").Value = "0.06710"
$ws.Range("E45").Value = "  -2.87%  "

$ws.Range("D46").Value = "0.4855"
$ws.Range("E46").Value = "  -6.38%  "

$ws.Range("D47").Value = "10.43"
$ws.Range("E47").Value = "  -6.86%  "

$ws.Range("D48").Value = "106.57"
$ws.Range("E48").Value = "  -3.48%  "

$ws.Range("D49").Value = "1.004"
$ws.Range("E49").Value = "  +0.04%  "

$ws.Range("D50").Value = "1.666"
$ws.Range("E50").Value = "  -5.67%  "

$ws.Range("D51").Value = "1.867"
$ws.Range("E51").Value = "  -15.29%  "
